$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended to the bottom of column A (rows 14-24)
$ws.Range("A14").Value = "PyShortTextCategorization"
$ws.Range("A15").Value = "image_super_resolution"
$ws.Range("A16").Value = "Flowise"
$ws.Range("A17").Value = "ChatVRM"
$ws.Range("A18").Value = "Umami"
$ws.Range("A19").Value = "ElasticSearch"
$ws.Range("A20").Value = "ReactNative"
$ws.Range("A21").Value = "SmartTube"
$ws.Range("A22").Value = "QuixBugs(Java)"
$ws.Range("A23").Value = "javascript-algorithms"

# Stray backtick value in M15 (added before the last row so shared-string
# order matches: ... javascript-algorithms, `, javascript-algorithms-2)
$ws.Range("M15").Value = "``"

$ws.Range("A24").Value = "javascript-algorithms-2"

# Column A width (ColumnWidth property is offset from the stored XML width
# by ~5/6 of a character for the default Calibri 11 font; 28.1666667 -> 29)
$ws.Columns.Item(1).ColumnWidth = 28.166666666666668

# Selection
$ws.Range("I20").Select()
